$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for Inhba-Acvr1 LR-pair sheet (YoungD7)
$updates = @{
    "E2" = 2
    "F2" = 0.6666666666666666
    "G2" = 0.1285113333333333
    "H2" = 0.385534
    "I2" = 0.03749201237720504
    "J2" = 0.03749201237720504
    "M2" = 10.92359866666667
    "N2" = 32.770796
    "O2" = 0.2236009040380497
    "P2" = 0.2236009040380497
    "Q2" = 1.403806229451556
    "R2" = 12.634256065064
    "S2" = 0.008383247861748794
    "T2" = 0.008383247861748796
    "E3" = 2
    "F3" = 0.6666666666666666
    "G3" = 0.1285113333333333
    "H3" = 0.385534
    "I3" = 0.03749201237720504
    "J3" = 0.03749201237720504
    "O3" = 0.4261214970992155
    "P3" = 0.4261214970992155
    "Q3" = 2.675266518731556
    "R3" = 24.077398668584
    "S3" = 0.01597615244343693
    "T3" = 0.01597615244343693
    "E4" = 2
    "F4" = 0.6666666666666666
    "G4" = 0.1285113333333333
    "H4" = 0.385534
    "I4" = 0.03749201237720504
    "J4" = 0.03749201237720504
    "M4" = 13.06524766666667
    "N4" = 39.195743
    "O4" = 0.2674394472823625
    "P4" = 0.2674394472823625
    "Q4" = 1.679032397973556
    "R4" = 15.111291581762
    "S4" = 0.01002684306766321
    "T4" = 0.01002684306766321
    "E5" = 2
    "F5" = 0.6666666666666666
    "G5" = 0.1285113333333333
    "H5" = 0.385534
    "I5" = 0.03749201237720504
    "J5" = 0.03749201237720504
    "M5" = 4.046901
    "N5" = 12.140703
    "O5" = 0.0828381515803724
    "P5" = 0.0828381515803724
    "Q5" = 0.5200726433780001
    "R5" = 4.680653790402
    "S5" = 0.003105769004356109
    "T5" = 0.00310576900435611
    "I6" = 0.7552862722193517
    "J6" = 0.755286272219352
    "M6" = 10.92359866666667
    "N6" = 32.770796
    "O6" = 0.2236009040380497
    "P6" = 0.2236009040380497
    "Q6" = 28.28003904654134
    "R6" = 254.520351418872
    "S6" = 0.1688826932757755
    "T6" = 0.1688826932757756
    "I7" = 0.7552862722193517
    "J7" = 0.755286272219352
    "O7" = 0.4261214970992155
    "P7" = 0.4261214970992155
    "S7" = 0.3218437170565958
    "T7" = 0.3218437170565959
    "I8" = 0.7552862722193517
    "J8" = 0.755286272219352
    "M8" = 13.06524766666667
    "N8" = 39.195743
    "O8" = 0.2674394472823625
    "P8" = 0.2674394472823625
    "Q8" = 33.82454129274733
    "R8" = 304.420871634726
    "S8" = 0.2019933431822994
    "T8" = 0.2019933431822995
    "I9" = 0.7552862722193517
    "J9" = 0.755286272219352
    "M9" = 4.046901
    "N9" = 12.140703
    "O9" = 0.0828381515803724
    "P9" = 0.0828381515803724
    "Q9" = 10.476997717494
    "R9" = 94.29297945744599
    "S9" = 0.06256651870468108
    "T9" = 0.06256651870468109
    "G10" = 0.692415
    "H10" = 2.077245
    "I10" = 0.2020057770533527
    "J10" = 0.2020057770533527
    "M10" = 10.92359866666667
    "N10" = 32.770796
    "O10" = 0.2236009040380497
    "P10" = 0.2236009040380497
    "Q10" = 7.563663570780001
    "R10" = 68.07297213702
    "S10" = 0.04516867437003838
    "T10" = 0.04516867437003838
    "G11" = 0.692415
    "H11" = 2.077245
    "I11" = 0.2020057770533527
    "J11" = 0.2020057770533527
    "O11" = 0.4261214970992155
    "P11" = 0.4261214970992155
    "Q11" = 14.41425140118
    "R11" = 129.72826261062
    "S11" = 0.08607900414066502
    "T11" = 0.08607900414066502
    "G12" = 0.692415
    "H12" = 2.077245
    "I12" = 0.2020057770533527
    "J12" = 0.2020057770533527
    "M12" = 13.06524766666667
    "N12" = 39.195743
    "O12" = 0.2674394472823625
    "P12" = 0.2674394472823625
    "Q12" = 9.046573463114999
    "R12" = 81.419161168035
    "S12" = 0.0540243133629928
    "T12" = 0.05402431336299281
    "G13" = 0.692415
    "H13" = 2.077245
    "I13" = 0.2020057770533527
    "J13" = 0.2020057770533527
    "M13" = 4.046901
    "N13" = 12.140703
    "O13" = 0.0828381515803724
    "P13" = 0.0828381515803724
    "Q13" = 2.802134955915
    "R13" = 25.219214603235
    "S13" = 0.01673378517965654
    "T13" = 0.01673378517965655
    "E14" = 1
    "F14" = 0.3333333333333333
    "G14" = 0.01787866666666667
    "H14" = 0.053636
    "I14" = 0.005215938350090445
    "J14" = 0.005215938350090446
    "M14" = 10.92359866666667
    "N14" = 32.770796
    "O14" = 0.2236009040380497
    "P14" = 0.2236009040380497
    "Q14" = 0.1952993793617778
    "R14" = 1.757694414256
    "S14" = 0.001166288530486957
    "T14" = 0.001166288530486957
    "E15" = 1
    "F15" = 0.3333333333333333
    "G15" = 0.01787866666666667
    "H15" = 0.053636
    "I15" = 0.005215938350090445
    "J15" = 0.005215938350090446
    "O15" = 0.4261214970992155
    "P15" = 0.4261214970992155
    "Q15" = 0.3721866164817778
    "R15" = 3.349679548336
    "S15" = 0.002222623458517753
    "T15" = 0.002222623458517753
    "E16" = 1
    "F16" = 0.3333333333333333
    "G16" = 0.01787866666666667
    "H16" = 0.053636
    "I16" = 0.005215938350090445
    "J16" = 0.005215938350090446
    "M16" = 13.06524766666667
    "N16" = 39.195743
    "O16" = 0.2674394472823625
    "P16" = 0.2674394472823625
    "Q16" = 0.2335892079497778
    "R16" = 2.102302871548
    "S16" = 0.001394947669407067
    "T16" = 0.001394947669407067
    "E17" = 1
    "F17" = 0.3333333333333333
    "G17" = 0.01787866666666667
    "H17" = 0.053636
    "I17" = 0.005215938350090445
    "J17" = 0.005215938350090446
    "M17" = 4.046901
    "N17" = 12.140703
    "O17" = 0.0828381515803724
    "P17" = 0.0828381515803724
    "Q17" = 0.07235319401200001
    "R17" = 0.6511787461080001
    "S17" = 0.0004320786916786698
    "T17" = 0.0004320786916786699
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

